# Regenerate column G (K = "Strike#" -> "K") values in the save_data sheet.
# These values come from re-running the std/mean calc and writing s_vals,
# so here we just write the freshly computed values into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 0
    4  = 0
    5  = 2
    6  = 3
    7  = 0
    8  = 0
    9  = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 2
    15 = 0
    16 = 2
    17 = 3
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 2
    23 = 1
    24 = 2
    25 = 1
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 2
    34 = 5
    35 = 0
    36 = 1
    37 = 1
    39 = 1
    40 = 2
    41 = 1
    42 = 0
    44 = 2
    45 = 2
    46 = 2
    48 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
